{"js": "// Update the \"TH\u1ed0NG K\u00ca COI THI K\u1ebeT TH\u00daC H\u1eccC PH\u1ea6N\" stats table:\n//  1. Remove the \"L\u1edbp 10A1\" data row entirely.\n//  2. Update the \"L\u1edbp 10A2\" row's exam date from 27-12-2018 to 07-01-2019.\n//  3. Insert a new \"L\u1edbp 11A1\" data row right after the \"L\u1edbp 10A2\" row.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The second table (index 1) is the statistics grid with the per-class rows.\nconst table = tables.items[1];\n\nasync function loadRowTexts() {\n  table.rows.load(\"items\");\n  await context.sync();\n  const rows = table.rows.items;\n  for (let i = 0; i < rows.length; i++) {\n    rows[i].cells.load(\"items\");\n  }\n  await context.sync();\n  for (let i = 0; i < rows.length; i++) {\n    const cells = rows[i].cells.items;\n    for (let j = 0; j < cells.length; j++) {\n      cells[j].body.load(\"text\");\n    }\n  }\n  await context.sync();\n  return rows;\n}\n\nfunction firstCellText(row) {\n  const cells = row.cells.items;\n  for (let j = 0; j < cells.length; j++) {\n    const text = cells[j].body.text.trim();\n    if (text.length > 0) return text;\n  }\n  return \"\";\n}\n\n// 1. Find and delete the \"L\u1edbp 10A1\" row.\nlet rows = await loadRowTexts();\nlet idx10A1 = -1;\nfor (let i = 0; i < rows.length; i++) {\n  if (firstCellText(rows[i]) === \"L\u1edbp 10A1\") idx10A1 = i;\n}\nif (idx10A1 >= 0) {\n  rows[idx10A1].delete();\n  await context.sync();\n}\n\n// Re-fetch the row collection: after the delete, previously-held row\n// references point at stale table positions.\nrows = await loadRowTexts();\nlet idx10A2 = -1;\nfor (let i = 0; i < rows.length; i++) {\n  if (firstCellText(rows[i]) === \"L\u1edbp 10A2\") idx10A2 = i;\n}\n\nif (idx10A2 >= 0) {\n  const row10A2 = rows[idx10A2];\n\n  // 2. Update the exam date (4th cell: TT | L\u1edbp | Qu\u00e2n s\u1ed1 | Ng\u00e0y thi | ...).\n  row10A2.cells.items[3].value = \"07-01-2019\";\n  await context.sync();\n\n  // 3. Add the new \"L\u1edbp 11A1\" row right after the \"L\u1edbp 10A2\" row.\n  row10A2.insertRows(\"After\", 1, [\n    [\"\", \"L\u1edbp 11A1\", \"4\", \"10-01-2019\", \"To\u00e1n\", \"Tr\u1eafc nghi\u1ec7m\", \"Dungvv\", \"T\u0103ng Thi\u1ebft Gi\u00e1p\", \"Th\u01b0\u1ee3ng \u00fay\", \" \"]\n  ]);\n  await context.sync();\n}\n", "ps1": "# Update the \"TH\u1ed0NG K\u00ca COI THI K\u1ebeT TH\u00daC H\u1eccC PH\u1ea6N\" stats table:\n#  1. Remove the \"L\u1edbp 10A1\" data row entirely.\n#  2. Update the \"L\u1edbp 10A2\" row's exam date from 27-12-2018 to 07-01-2019.\n#  3. Insert a new \"L\u1edbp 11A1\" data row right after the \"L\u1edbp 10A2\" row.\n\n$d = $word.ActiveDocument\n\n# The second table is the statistics grid with the per-class rows.\n$t = $d.Tables.Item(2)\n\nfunction Get-RowLabel($row) {\n    for ($c = 1; $c -le $row.Cells.Count; $c++) {\n        $txt = $row.Cells.Item($c).Range.Text\n        $txt = $txt -replace \"[\\r\\a\\x07]\", \"\"\n        $txt = $txt.Trim()\n        if ($txt.Length -gt 0) {\n            return $txt\n        }\n    }\n    return \"\"\n}\n\n# 1. Find and delete the \"L\u1edbp 10A1\" row.\n$idx10A1 = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    if ((Get-RowLabel($t.Rows.Item($r))) -eq \"L\u1edbp 10A1\") {\n        $idx10A1 = $r\n    }\n}\nif ($idx10A1 -gt 0) {\n    $t.Rows.Item($idx10A1).Delete()\n}\n\n# 2. Find the \"L\u1edbp 10A2\" row and update its exam-date cell (4th cell:\n#    TT | L\u1edbp | Qu\u00e2n s\u1ed1 | Ng\u00e0y thi | ...).\n$idx10A2 = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    if ((Get-RowLabel($t.Rows.Item($r))) -eq \"L\u1edbp 10A2\") {\n        $idx10A2 = $r\n    }\n}\n\nif ($idx10A2 -gt 0) {\n    $row10A2 = $t.Rows.Item($idx10A2)\n    $row10A2.Cells.Item(4).Range.Text = \"07-01-2019\"\n\n    # 3. Insert the new \"L\u1edbp 11A1\" row right after the \"L\u1edbp 10A2\" row. If\n    #    \"L\u1edbp 10A2\" is the last row, Rows.Add() with no argument appends a\n    #    fresh row at the very end, which is exactly \"right after\" it;\n    #    otherwise insert before the row that currently follows it.\n    if ($idx10A2 -eq $t.Rows.Count) {\n        $newRow = $t.Rows.Add()\n    } else {\n        $newRow = $t.Rows.Add($t.Rows.Item($idx10A2 + 1))\n    }\n\n    $values = @(\"\", \"L\u1edbp 11A1\", \"4\", \"10-01-2019\", \"To\u00e1n\", \"Tr\u1eafc nghi\u1ec7m\", \"Dungvv\", \"T\u0103ng Thi\u1ebft Gi\u00e1p\", \"Th\u01b0\u1ee3ng \u00fay\", \" \")\n    for ($c = 1; $c -le $values.Length; $c++) {\n        $newRow.Cells.Item($c).Range.Text = $values[$c - 1]\n    }\n}\n"}
